# "leading_quotes" test sheet refresh:
#  - shrink the LEFT() demo range from C4:C21 down to C4:C18
#  - add a descriptive label column (A) explaining each B-column example
#  - turn the blank quote-prefixed B5 example into a real "def" value
#  - relocate the stray quote-prefixed blank formatting demo from F12 to F9
#  - nudge column widths/selection to match the refreshed layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 19-21 (bottom of the shared formula range)
$ws.Rows("19:21").Delete()

# Set new values for column A (order matters for shared-string table index)
$ws.Range("A4").Value = "Quote prepended string"
$ws.Range("A3").Value = "Quote prepended string that looks like a number"
$ws.Range("A5").Value = "String without quote"
$ws.Range("B5").Value = "def"

# Adjust column widths
$ws.Columns("A").ColumnWidth = 39.83
$ws.Columns("B").ColumnWidth = 11.0
$ws.Columns("C").ColumnWidth = 12.67

# Move the quote-prefixed empty cell style from F12 to F9
$ws.Range("F12").Copy()
$ws.Range("F9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F12").Clear()

# Update selection
[void]$ws.Range("C4").Select()

Write-Host "done"
